# Generate Report for Handoff
#
# Updates the localization-status report so the Overview/zh-cn/de-de sheets
# reflect that the content is ready for handoff (rather than "Handed back:
# in sync with en-US"), and refreshes the associated generation/handoff
# timestamps.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-31 09:15:00"

# Columns narrow now that the status text is shorter.
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-31 09:14:55"
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-31 09:15:00"
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
